$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (primera propuesta) ---
# Cliente
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "23162539"

# Monto
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "4500"

# Numero Cuotas
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "24"

# Numero Propuesta (resultado de la emision del dictamen)
$ws.Range("V2").NumberFormat = "@"
$ws.Range("V2").Value = "4900002"
$ws.Range("V2").Style = "Normal"

# Resultado
$ws.Range("W2").NumberFormat = "@"
$ws.Range("W2").Value = " "
$ws.Range("W2").Style = "Normal"

# --- Row 3 (segunda propuesta) ---
# Cliente
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "12668309"

# Monto
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "1500"

# Numero Cuotas
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "12"

# Numero Propuesta / Resultado: todavia no se ha emitido dictamen, quedan vacios
$ws.Range("V3").ClearContents()
$ws.Range("W3").ClearContents()

# Restablecer la vista de la hoja: sin desplazamiento y con el cursor en A1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A1").Select() | Out-Null
